$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.79"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.31%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.54%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.039"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.16%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07849"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.32%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.205"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.49%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.992"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.79%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9263"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.30%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09849"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.30%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1886"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.16%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08662"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.91%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03684"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "8.75%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09937"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.05%"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.61%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005648"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.61%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.462"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.64%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.019"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.03%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.253"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "7.54%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3411"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.53%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1304"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.51%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.766"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.93%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2207"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.48%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04602"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.48%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001257"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.19%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004479"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.34%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001403"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "8.10%"

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002723"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-19.72%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01844"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.64%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04759"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.70%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.008042"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.98%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1400"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.49%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007575"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-13.87%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002102"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.10%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01041"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "13.48%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006284"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "3.66%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.41%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005818"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.30%"

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.21"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "882.77%"

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.002695"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.37%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002105"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.41%"
